# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" sheet and the corresponding rows on the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1597
$ws1.Range("F4").Value = 5221
$ws1.Range("F6").Value = 10390
$ws1.Range("F7").Value = 266
$ws1.Range("F8").Value = 566
$ws1.Range("F9").Value = 119
$ws1.Range("F10").Value = 128
$ws1.Range("F11").Value = 825

# --- Sheet "全部类型" (all types, combined rows) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1597
$ws4.Range("F6").Value = 5221
$ws4.Range("F9").Value = 10390
$ws4.Range("F10").Value = 266
$ws4.Range("F11").Value = 566
$ws4.Range("F12").Value = 119
$ws4.Range("F15").Value = 128
$ws4.Range("F16").Value = 825
